# "Hjemme passive updated meanEMG legmaxROM"
# The sheet originally held two repeated 31-column blocks of "degree" headers
# (cols B..Z and AA..AY). This edit collapses the view down to just the
# "15"/"16" degree columns (taken from the duplicate block), reusing columns
# B:E, and updates the corresponding CON/STR rows with the matching values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (header / degree labels)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 (CON) - B2/D2 updated, C2/E2 removed entirely
$ws.Range("B2").Value = 20.397789610916238
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 23.26160720194569
$ws.Range("E2").ClearContents()

# Row 3 (STR) - all four values updated
$ws.Range("B3").Value = 16.455616307195758
$ws.Range("C3").Value = -7.741317205820601
$ws.Range("D3").Value = 14.833894731861719
$ws.Range("E3").Value = -13.611244854621201

# Selection now only spans the trimmed range
$ws.Range("B1:E3").Select()
